$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: selection changed to column Q ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("Q").Select()

# --- Transactions sheet: data + formatting updates ---
$ws = $wb.Worksheets.Item("Transactions")

# Update transaction IDs
$ws.Range("A2").Value = 1454
$ws.Range("A3").Value = 1453
$ws.Range("A4").Value = 1449
$ws.Range("A5").Value = 1448

# Update interest amounts on row 2
$ws.Range("E2").Value = 7.96
$ws.Range("I2").Value = 7.96

# Apply italic formatting (matching existing data style: vertical-center + wrap)
# to the new K:L columns on rows 2-4, copying the base formatting from J2 first
# so the resulting style keeps the same alignment as the rest of the row.
$ws.Range("J2").Copy()
$ws.Range("K2:L4").PasteSpecial(-4122)
$ws.Range("K2:L4").Font.Italic = $true

# Touch K5:L5 so they exist as empty cells with the default (Normal) style
$ws.Range("K5:L5").Style = "Normal"

# Update the active selection on the Transactions sheet
$ws.Range("I5").Select()

$excel.CutCopyMode = $false
